$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 160251
$ws.Range("C4").Value = 151287
$ws.Range("C5").Value = 8965
$ws.Range("C7").Value = 5.59
$ws.Range("C8").Value = 64.40000000000001
